$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1) — same bold/centered/bordered style as existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header formatting (bold font, border, center/top alignment) from an
# existing header cell onto the new ones without touching their values.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier-flag data for rows 2-12, columns F (KNN), G (SVM), H (RF)
$values = @(
    @(2,  $false, $false, $false),
    @(3,  $false, $false, $false),
    @(4,  $false, $false, $false),
    @(5,  $false, $false, $false),
    @(6,  $true,  $true,  $false),
    @(7,  $false, $false, $false),
    @(8,  $false, $false, $false),
    @(9,  $false, $false, $false),
    @(10, $true,  $false, $false),
    @(11, $false, $false, $false),
    @(12, $false, $false, $false)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
}
